# --- Update Rushing sheet (Week 15 stat logging) ---
$wb = $excel.ActiveWorkbook
$rushing = $wb.Worksheets.Item("Rushing")
$rushing.Range("C2").Value = 6
$rushing.Range("D2").Value = 5
$rushing.Range("E2").Value = 9
$rushing.Range("F2").Value = 6
$rushing.Range("C3").Value = 9
$rushing.Range("D3").Value = 6
$rushing.Range("E3").Value = 5
$rushing.Range("C6").Value = 50
$rushing.Range("D6").Value = 31
$rushing.Range("E6").Value = 17
$rushing.Range("F6").Value = 21
$rushing.Range("C7").Value = 28
$rushing.Range("D7").Value = 17
$rushing.Range("F7").Value = 9

# --- Update Receiving sheet (Week 15 stat logging) ---
$receiving = $wb.Worksheets.Item("Receiving")
$receiving.Range("C2").Value = 14
$receiving.Range("D2").Value = 9
$receiving.Range("G2").Value = 4
$receiving.Range("H2").Value = 3
$receiving.Range("C4").Value = 8
$receiving.Range("D4").Value = 7
$receiving.Range("C5").Value = 20
$receiving.Range("D5").Value = 13
$receiving.Range("G5").Value = 5
$receiving.Range("H5").Value = 4
$receiving.Range("C6").Value = 49
$receiving.Range("G6").Value = 6
$receiving.Range("C7").Value = 26
$receiving.Range("D7").Value = 17
$receiving.Range("C8").Value = 2
$receiving.Range("D8").Value = 2
$receiving.Range("C11").Value = 22
$receiving.Range("D11").Value = 18
$receiving.Range("E11").Value = 6
$receiving.Range("G11").Value = 4
$receiving.Range("H11").Value = 4
$receiving.Range("C12").Value = 6
$receiving.Range("D12").Value = 4
$receiving.Range("E12").Value = 1
$receiving.Range("C13").Value = 28
$receiving.Range("D13").Value = 23
$receiving.Range("G13").Value = 2
$receiving.Range("H13").Value = 1
$receiving.Range("C14").Value = 14
$receiving.Range("D14").Value = 10
$receiving.Range("C15").Value = 26
$receiving.Range("D15").Value = 20
$receiving.Range("G15").Value = 6

# --- Add new player row 17 (Week 16 simulated: A.Brown) ---
$receiving.Range("A17").Value = 15
$receiving.Range("B17").Value = "A.Brown"
$receiving.Range("C17").Value = 80
$receiving.Range("D17").Value = 60
$receiving.Range("E17").Value = 18
$receiving.Range("F17").Value = 15
$receiving.Range("G17").Value = 9
$receiving.Range("H17").Value = 8

# --- Style the new row number cell like the others (bold, centered, bordered) ---
$a17 = $receiving.Range("A17")
$a17.Font.Bold = $true
$a17.HorizontalAlignment = -4108
$a17.VerticalAlignment = -4160
$a17.Borders.Item(7).LineStyle = 1
$a17.Borders.Item(10).LineStyle = 1

# --- Switch active tab / selection to Receiving sheet ---
$receiving.Activate()
$receiving.Range("C21").Select()
